# Add the "Student Body Makeup" worksheet (student demographics) after the
# existing "Aggregated Enrollment Data" sheet, and populate it with the new
# geographic / gender / ethnic breakdown data.

$wb = $excel.ActiveWorkbook
$wsAgg = $wb.Worksheets.Item(1)

$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$ws.Name = "Student Body Makeup"

# --- Seed the shared-string table in the same order the source data was
#     authored, so that cells sharing a string (Option/Data labels, notes)
#     resolve to the same underlying entries. ------------------------------

$ws.Range("A1").Value = "Data"

$ws.Range("B3").Value = "Other U.S."
$ws.Range("B2").Value = "California"
$ws.Range("B4").Value = "International"
$ws.Range("A2").Value = "Geographic Origin (%)"

$ws.Range("B11").Value = "Female"
$ws.Range("A11").Value = "Gender Balance (%)"

$ws.Range("A13").Value = "Ethnic Diversity (%)"
$ws.Range("B20").Value = "Two or more"
$ws.Range("B13").Value = "African American"
$ws.Range("B14").Value = "Asian"
$ws.Range("B16").Value = "Native American"
$ws.Range("B17").Value = "Native Hawaiian/Pacific Islander"
$ws.Range("B19").Value = "White"
$ws.Range("B21").Value = "Unknown"

$ws.Range("J8").Value = "Asian/Asian American was used interchangably by Stanford. A distiction cannot be made here."
$ws.Range("J15").Value = "Also listed as ""Declined to State/Other"" in some years"
$ws.Range("J12").Value = "Also listed as ""Hispanic/Latino"" in some years"
$ws.Range("B18").Value = "Hispanic"
$ws.Range("J3").Value = "Stateless peoples are counted as ""Other U.S."""

$ws.Range("B6").Value = "The Americas"
$ws.Range("B7").Value = "Europe"
$ws.Range("B8").Value = "Africa"
$ws.Range("B9").Value = "Middle East and North Africa"
$ws.Range("B10").Value = "Pacific Basin"
$ws.Range("B5").Value = "Asia"
$ws.Range("A5").Value = "Geographic Origin (% of Foreign)"

$ws.Range("J2").Value = "Measured in the Fall of each year, describes entire undergraduate class, data taken from https://facts.stanford.edu/academics/undergraduate-profile/ with past data retreived via the WayBack Machine (eg. https://web.archive.org/web/20130406141549/http://facts.stanford.edu/academics/undergraduate-profile)"

# --- Fill in the rest of the labels that reuse already-created strings ----

$ws.Range("A3").Value = "Geographic Origin (%)"
$ws.Range("A4").Value = "Geographic Origin (%)"
$ws.Range("A6").Value = "Geographic Origin (% of Foreign)"
$ws.Range("A7").Value = "Geographic Origin (% of Foreign)"
$ws.Range("A8").Value = "Geographic Origin (% of Foreign)"
$ws.Range("A9").Value = "Geographic Origin (% of Foreign)"
$ws.Range("A10").Value = "Geographic Origin (% of Foreign)"
$ws.Range("A12").Value = "Gender Balance (%)"
$ws.Range("B12").Value = "Men"
$ws.Range("A14").Value = "Ethnic Diversity (%)"
$ws.Range("A15").Value = "Ethnic Diversity (%)"
$ws.Range("B15").Value = "International"
$ws.Range("A16").Value = "Ethnic Diversity (%)"
$ws.Range("A17").Value = "Ethnic Diversity (%)"
$ws.Range("A18").Value = "Ethnic Diversity (%)"
$ws.Range("A19").Value = "Ethnic Diversity (%)"
$ws.Range("A20").Value = "Ethnic Diversity (%)"
$ws.Range("A21").Value = "Ethnic Diversity (%)"

# --- Header row -------------------------------------------------------------

$ws.Range("B1").Value = "Option"
$ws.Range("C1").Value = "2012-2013"
$ws.Range("D1").Value = "2013-2014"
$ws.Range("E1").Value = "2014-2015"
$ws.Range("F1").Value = "2015-2016"
$ws.Range("G1").Value = "2016-2017"
$ws.Range("H1").Value = "2017-2018"
$ws.Range("I1").Value = "2018-2019"
$ws.Range("J1").Value = "Notes"

# --- Row 2: Geographic Origin (%) / California -----------------------------
$ws.Range("C2").Value = 40
$ws.Range("D2").Value = 39
$ws.Range("E2").Value = 37.6
$ws.Range("F2").Value = 36
$ws.Range("G2").Value = 36
$ws.Range("H2").Value = 35
$ws.Range("I2").Value = 35

# --- Row 3: Geographic Origin (%) / Other U.S. ------------------------------
$ws.Range("C3").Value = 53
$ws.Range("D3").Value = 53
$ws.Range("E3").Value = 54.2
$ws.Range("F3").Value = 55.3
$ws.Range("G3").Value = 53
$ws.Range("H3").Value = 53
$ws.Range("I3").Value = 52

# --- Row 4: Geographic Origin (%) / International ---------------------------
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 8.2
$ws.Range("F4").Value = 8.7
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 12
$ws.Range("I4").Value = 13

# --- Row 5: Geographic Origin (% of Foreign) / Asia -------------------------
$ws.Range("C5").Value = 52
$ws.Range("D5").Value = 50
$ws.Range("E5").Value = 47.6
$ws.Range("F5").Value = 43.79
$ws.Range("G5").Value = "N/G"
$ws.Range("H5").Value = "N/G"
$ws.Range("I5").Value = "N/G"

# --- Row 6: Geographic Origin (% of Foreign) / The Americas -----------------
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 19
$ws.Range("E6").Value = 20.4
$ws.Range("F6").Value = 22.71
$ws.Range("G6").Value = "N/G"
$ws.Range("H6").Value = "N/G"
$ws.Range("I6").Value = "N/G"

# --- Row 7: Geographic Origin (% of Foreign) / Europe -----------------------
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 16
$ws.Range("E7").Value = 16.4
$ws.Range("F7").Value = 16.34
$ws.Range("G7").Value = "N/G"
$ws.Range("H7").Value = "N/G"
$ws.Range("I7").Value = "N/G"

# --- Row 8: Geographic Origin (% of Foreign) / Africa -----------------------
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = 7.8
$ws.Range("F8").Value = 8.5
$ws.Range("G8").Value = "N/G"
$ws.Range("H8").Value = "N/G"
$ws.Range("I8").Value = "N/G"

# --- Row 9: Geographic Origin (% of Foreign) / Middle East and North Africa -
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 5.2
$ws.Range("F9").Value = 6.05
$ws.Range("G9").Value = "N/G"
$ws.Range("H9").Value = "N/G"
$ws.Range("I9").Value = "N/G"

# --- Row 10: Geographic Origin (% of Foreign) / Pacific Basin ---------------
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 2.6
$ws.Range("F10").Value = 2.61
$ws.Range("G10").Value = "N/G"
$ws.Range("H10").Value = "N/G"
$ws.Range("I10").Value = "N/G"

# --- Row 11: Gender Balance (%) / Female ------------------------------------
$ws.Range("C11").Value = 48
$ws.Range("D11").Value = 47
$ws.Range("E11").Value = 47.2
$ws.Range("F11").Value = 47.6
$ws.Range("G11").Value = 48
$ws.Range("H11").Value = 50
$ws.Range("I11").Value = 50

# --- Row 12: Gender Balance (%) / Men ---------------------------------------
$ws.Range("C12").Value = 52
$ws.Range("D12").Value = 53
$ws.Range("E12").Value = 52.8
$ws.Range("F12").Value = 52.4
$ws.Range("G12").Value = 52
$ws.Range("H12").Value = 50
$ws.Range("I12").Value = 50

# --- Row 13: Ethnic Diversity (%) / African American ------------------------
$ws.Range("C13").Value = 9
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = 7.5
$ws.Range("F13").Value = 7.8
$ws.Range("G13").Value = 6
$ws.Range("H13").Value = 7
$ws.Range("I13").Value = 6

# --- Row 14: Ethnic Diversity (%) / Asian -----------------------------------
$ws.Range("C14").Value = 21
$ws.Range("D14").Value = 22
$ws.Range("E14").Value = 22.6
$ws.Range("F14").Value = 22.9
$ws.Range("G14").Value = 21
$ws.Range("H14").Value = 22
$ws.Range("I14").Value = 22

# --- Row 15: Ethnic Diversity (%) / International ---------------------------
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = 8.2
$ws.Range("F15").Value = 8.8
$ws.Range("G15").Value = 9
$ws.Range("H15").Value = 9
$ws.Range("I15").Value = 10

# --- Row 16: Ethnic Diversity (%) / Native American -------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 2.1
$ws.Range("F16").Value = 1.9
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1

# --- Row 17: Ethnic Diversity (%) / Native Hawaiian/Pacific Islander --------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1.2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.5
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 0.5

# --- Row 18: Ethnic Diversity (%) / Hispanic --------------------------------
$ws.Range("C18").Formula = "=6+7"
$ws.Range("D18").Formula = "=7+7"
$ws.Range("E18").Formula = "=6.8+6.3"
$ws.Range("F18").Value = 12.6
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 15
$ws.Range("I18").Value = 16

# --- Row 19: Ethnic Diversity (%) / White -----------------------------------
$ws.Range("C19").Value = 39
$ws.Range("D19").Value = 41
$ws.Range("E19").Value = 42.8
$ws.Range("F19").Value = 42.5
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 36
$ws.Range("I19").Value = 65

# --- Row 20: Ethnic Diversity (%) / Two or more -----------------------------
$ws.Range("C20").Value = "N/G"
$ws.Range("D20").Value = "N/G"
$ws.Range("E20").Value = "N/G"
$ws.Range("F20").Value = "N/G"
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 9
$ws.Range("I20").Value = 9

# --- Row 21: Ethnic Diversity (%) / Unknown ---------------------------------
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 2.6
$ws.Range("F21").Value = 2.5
$ws.Range("G21").Value = 0.5
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 0.5

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.451822916666668
$ws.Columns.Item(2).ColumnWidth = 26.877604166666668

# --- Selection / active-sheet bookkeeping -----------------------------------
# "Aggregated Enrollment Data" keeps a plain selection (no longer the active tab)
$wsAgg.Range("G1:M1").Select()
# The new sheet is the one left active/selected when the workbook is saved.
$ws.Range("A6").Select()
